# Update to new files:
#   - The two "F1" rows are disambiguated into "F1a" and "F1b"
#   - The two "F2" rows are disambiguated into "F2a" and "F2b"
#   - Selection moves from H12 to C16

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the cross labels; rows 4-5 were both "F1" and rows 6-7
# were both "F2". Replace them with distinguishable labels.
$ws.Range("A4").Value = "F1a"
$ws.Range("A5").Value = "F1b"
$ws.Range("A6").Value = "F2a"
$ws.Range("A7").Value = "F2b"

# Update the active selection to match the saved view.
$ws.Range("C16").Select()
